# Apply odds/value updates to Sheet1 as described by the commit diff
# (71 individual cell value changes across rows 2, 3, 4, 5, 6, 8, 9, 10)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("J2").Value = 1.95
$ws.Range("Q2").Value = 1.53
$ws.Range("R2").Value = 2.5

# Row 3
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 3
$ws.Range("Q3").Value = 2.2
$ws.Range("R3").Value = 1.67
$ws.Range("AT3").Value = 2.63

# Row 4
$ws.Range("G4").Value = 1.85
$ws.Range("H4").Value = 3.25
$ws.Range("I4").Value = 4.5
$ws.Range("J4").Value = 2.63
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 5
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 8
$ws.Range("U4").Value = 2.2
$ws.Range("V4").Value = 1.62
$ws.Range("W4").Value = 5.5
$ws.Range("X4").Value = 7.5
$ws.Range("Y4").Value = 9
$ws.Range("Z4").Value = 15
$ws.Range("AB4").Value = 34
$ws.Range("AE4").Value = 21
$ws.Range("AF4").Value = 81
$ws.Range("AH4").Value = 9.5
$ws.Range("AI4").Value = 21
$ws.Range("AK4").Value = 51
$ws.Range("AN4").Value = 3.6
$ws.Range("AO4").Value = 11
$ws.Range("AU4").Value = 9.5
$ws.Range("AV4").Value = 81
$ws.Range("AX4").Value = 6
$ws.Range("AY4").Value = 26
$ws.Range("BA4").Value = 101
$ws.Range("BB4").Value = 151

# Row 5
$ws.Range("G5").Value = 1.95
$ws.Range("H5").Value = 2.88
$ws.Range("J5").Value = 2.87
$ws.Range("K5").Value = 1.8
$ws.Range("L5").Value = 5
$ws.Range("O5").Value = 1.67
$ws.Range("P5").Value = 2.1
$ws.Range("U5").Value = 2.5
$ws.Range("V5").Value = 1.5
$ws.Range("W5").Value = 5
$ws.Range("X5").Value = 8
$ws.Range("AH5").Value = 7.5
$ws.Range("AI5").Value = 19
$ws.Range("AM5").Value = 51
$ws.Range("AO5").Value = 13
$ws.Range("AU5").Value = 10
$ws.Range("BA5").Value = 101
$ws.Range("BB5").Value = 151

# Row 6
$ws.Range("J6").Value = 1.87
$ws.Range("K6").Value = 2.62
$ws.Range("U6").Value = 1.62

# Row 8
$ws.Range("G8").Value = 1.38
$ws.Range("H8").Value = 6
$ws.Range("I8").Value = 5.75
$ws.Range("J8").Value = 1.77
$ws.Range("U8").Value = 1.36
$ws.Range("V8").Value = 3
$ws.Range("W8").Value = 19
$ws.Range("AF8").Value = 29
$ws.Range("AG8").Value = 67
$ws.Range("AJ8").Value = 19
$ws.Range("AY8").Value = 26
$ws.Range("BC8").Value = 81

# Row 9
$ws.Range("U9").Value = 1.57

# Row 10
$ws.Range("Q10").Value = 1.36

